$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "args" cell for the send_command row to reference the new
# s2_command() cmd_info string
$ws.Range("E5").Value = "s2_command() cmd_info"

# Autofit column E so its width reflects the new, longer text
$ws.Columns.Item(5).AutoFit() | Out-Null

$ws.Range("C24").Select() | Out-Null
